$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: "Testing" + " " + "custom" + " " + "properties"
#     -> "Testing " + "custom " + "properties"
$titleShape = $s.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

# Merge "Testing" and the following space into a single run "Testing "
$titleRange.Characters(1, 8).Text = "Testing "

# Merge "custom" and the following space into a single run "custom "
$titleShape.TextFrame.TextRange.Characters(9, 7).Text = "custom "

# --- Subtitle shape: "A." + " " + "M."  ->  "A. " + "M."
$subtitleShape = $s.Shapes.Item(2)
$subtitleRange = $subtitleShape.TextFrame.TextRange

# Merge "A." and the following space into a single run "A. "
$subtitleRange.Characters(3, 3).Text = "A. "
